# 12 6 21 rich mosla
# Add the Dec-12 "rich mosla" shopping entries (rows 9-13 of the
# December sheet) and move the selection down to F14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value2  = "2kg alu"
$ws.Range("F9").Value2  = 60

$ws.Range("E10").Value2 = "holud gura + Jira gura "
$ws.Range("F10").Value2 = 90

$ws.Range("E11").Value2 = "shak "
$ws.Range("F11").Value2 = 30

$ws.Range("E12").Value2 = "Chal"
$ws.Range("F12").Value2 = 260

$ws.Range("E13").Value2 = "Rishka"
$ws.Range("F13").Value2 = 20

$ws.Range("F14").Select()
